$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.300.36"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.57"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.46"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.12"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.506.31"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  -1.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  +8.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.587"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.13"
$ws.Range("E13").Value = "  -2.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000277"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.094.85"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "613.06"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.28"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.514.71"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.404.31"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.879"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("E23").Value = "  -8.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.33"
$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.61"
$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.56"
$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.17"
$ws.Range("E29").Value = "  +3.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.02"
$ws.Range("E30").Value = "  -1.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.05"
$ws.Range("E31").Value = "  -4.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.97"
$ws.Range("E32").Value = "  -2.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "658.90"
$ws.Range("E33").Value = "  +16.63%  "

$ws.Range("E34").Value = "  -4.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.81"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0996"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.74"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0476"
$ws.Range("E39").Value = "  +6.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.73"
$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("E42").Value = "  +1.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0740"
$ws.Range("E43").Value = "  +5.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.360.14"

$ws.Range("E45").Value = "  -4.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.07"
$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.55"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.66"
$ws.Range("E50").Value = "  -2.02%  "

$ws.Range("E51").Value = "  -0.01%  "
